$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F5").Value = 4799
$wsExhibition.Range("F10").Value = 216

# Sheet "演出" (performances)
$wsPerformance = $wb.Worksheets.Item("演出")
$wsPerformance.Range("F2").Value = 25
$wsPerformance.Range("F3").Value = 2

# Sheet "全部类型" (all types) - aggregated view of the above sheets
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 4799
$wsAll.Range("F10").Value = 25
$wsAll.Range("F11").Value = 216
$wsAll.Range("F12").Value = 2
